# Fix counting for room type for update flat type in officer.
# The "Number of units for Type 1" (2-Room) count for the first project
# ("Acacia Breeze") was wrong (0) and is corrected to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 1

# Leave the selection where the user last clicked after making the edit.
$ws.Range("F9").Select()
